$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = "newPrice"; E = "newVolume" }  (absent key = unchanged in that diff)
$updates = @{
    2 = @{ D="23.118.61"; E="  -3.58%  " }
    3 = @{ D="1.598.71"; E="  -3.27%  " }
    4 = @{ D="1.002"; E="  +0.34%  " }
    5 = @{ D="1.003"; E="  +0.38%  " }
    6 = @{ D="301.93"; E="  -2.28%  " }
    7 = @{ D="0.3772"; E="  -3.21%  " }
    8 = @{ D="0.3653"; E="  -4.69%  " }
    9 = @{ D="47.77"; E="  -6.74%  " }
    10 = @{ D="1.004"; E="  +0.53%  " }
    11 = @{ E="  -5.78%  " }
    12 = @{ D="0.08077"; E="  -4.34%  " }
    13 = @{ D="22.97"; E="  -4.00%  " }
    14 = @{ D="6.634"; E="  -6.88%  " }
    15 = @{ D="7.647"; E="  -2.46%  " }
    16 = @{ D="0.00001267"; E="  -3.91%  " }
    17 = @{ D="1.604.89"; E="  -2.76%  " }
    18 = @{ D="91.55"; E="  -3.12%  " }
    19 = @{ D="0.06794"; E="  -2.66%  " }
    20 = @{ D="18.42"; E="  -6.75%  " }
    21 = @{ D="6.589"; E="  -4.26%  " }
    22 = @{ E="  +0.25%  " }
    23 = @{ D="13.03"; E="  -4.18%  " }
    24 = @{ D="23.144.38"; E="  -3.42%  " }
    25 = @{ D="2.363"; E="  -4.76%  " }
    26 = @{ D="2.900"; E="  -3.85%  " }
    27 = @{ E="  -4.30%  " }
    28 = @{ D="151.06"; E="  -1.00%  " }
    29 = @{ D="5.247"; E="  -3.68%  " }
    30 = @{ D="131.83"; E="  -5.33%  " }
    31 = @{ D="2.443"; E="  -1.56%  " }
    32 = @{ D="7.141"; E="  -7.92%  " }
    33 = @{ D="1.771.84"; E="  -3.34%  " }
    34 = @{ D="0.9857"; E="  -4.67%  " }
    35 = @{ D="0.07731"; E="  -3.94%  " }
    36 = @{ D="0.02783"; E="  -6.23%  " }
    37 = @{ D="6.306"; E="  -5.91%  " }
    38 = @{ D="0.2547"; E="  -5.11%  " }
    39 = @{ D="0.08876"; E="  -2.64%  " }
    40 = @{ D="10.06"; E="  -7.24%  " }
    41 = @{ D="1.399"; E="  -1.95%  " }
    42 = @{ D="0.7169" }
    43 = @{ D="12.80"; E="  -5.03%  " }
    44 = @{ D="15.85"; E="  -2.04%  " }
    45 = @{ D="0.6646"; E="  -4.13%  " }
    46 = @{ D="2.313"; E="  -5.67%  " }
    47 = @{ E="  +0.26%  " }
    48 = @{ D="3.968"; E="  -2.60%  " }
    49 = @{ D="132.38"; E="  -1.51%  " }
    50 = @{ D="0.07977"; E="  -4.13%  " }
    51 = @{ D="1.174"; E="  -4.16%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        # Force text storage so numeric-looking prices ("1.002", "2.900") keep their
        # original literal digits instead of being parsed/normalised into a Number.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
        $cell.Style = "Normal"
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
